$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.579.48"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "3.316.88"
$ws.Range("E3").Value = "  +5.39%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.78"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.88"
$ws.Range("E6").Value = "  +2.86%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.316.46"
$ws.Range("E8").Value = "  +5.47%  "
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("E10").Value = "  +2.72%  "
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.67"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").Value = "3.868.08"
$ws.Range("E15").Value = "  +5.49%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "3.318.34"
$ws.Range("E17").Value = "  +5.51%  "
$ws.Range("D18").Value = "63.673.22"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.85"
$ws.Range("E19").Value = "  +2.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.65"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.13"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.736"
$ws.Range("E22").Value = "  +5.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.17"
$ws.Range("E23").Value = "  +5.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.70"
$ws.Range("E24").Value = "  +5.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.92"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +2.53%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.23"
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.17"
$ws.Range("E30").Value = "  +3.03%  "
$ws.Range("E31").Value = "  +2.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.90"
$ws.Range("E32").Value = "  +8.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.106"
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  +3.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.04"
$ws.Range("E36").Value = "  +4.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.66"
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("D38").Value = "0.0₃0742"
$ws.Range("E38").Value = "  +6.45%  "
$ws.Range("E39").Value = "  +2.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "433.97"
$ws.Range("E40").Value = "  +4.57%  "
$ws.Range("D41").Value = "3.090.54"
$ws.Range("E41").Value = "  +5.66%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.33"
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.117"
$ws.Range("E44").Value = "  +5.24%  "
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("E46").Value = "  +3.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.78"
$ws.Range("E47").Value = "  +14.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.29"
$ws.Range("E48").Value = "  +3.62%  "
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.84"
$ws.Range("E51").Value = "  +3.46%  "
